$d = $word.ActiveDocument

# Locate the paragraph that ends "...already assigned for other uses."
# via Find, then resolve it to a Paragraphs collection index so we can
# reliably manipulate whole paragraphs (style, insertion) afterwards.
$fr = $d.Content
$found = $fr.Find.Execute("already assigned for other uses.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetIndex = 0
if ($found) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $fr.Start -and $p.Range.End -ge $fr.End) {
            $targetIndex = $i
            break
        }
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $r = $target.Range
    $r.Collapse(0)

    # Insert two new empty paragraphs right after the target paragraph
    $r.InsertParagraphAfter()
    $r.Collapse(0)
    $r.InsertParagraphAfter()

    # The first new paragraph becomes the new Heading2
    $newHeading = $d.Paragraphs.Item($targetIndex + 1)
    $newHeading.Range.Text = "Window Activation on Switch From Empty Desktop"
    $newHeading.Style = "Heading2"

    # The second new paragraph becomes the body text (stays Normal, no explicit style)
    $newBody = $d.Paragraphs.Item($targetIndex + 2)
    $newBody.Range.Text = "Virtual Desktop Grid Switcher fixes an issue in Windows 10 where switching from a desktop which is empty to a desktop which had an activate window does not reactivate that window."
}
